$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.321.06"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "1.660.32"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").Value = "'219.82"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.65%  "
$ws.Range("D8").Value = "'0.256"
$ws.Range("E8").Value = "  +2.33%  "
$ws.Range("E9").Value = "  +0.50%  "
$ws.Range("D10").Value = "'20.06"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "1.891.76"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "1.661.68"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "'67.31"
$ws.Range("E16").Value = "  +4.54%  "
$ws.Range("D17").Value = "27.309.74"
$ws.Range("E17").Value = "  +2.14%  "
$ws.Range("D18").Value = "0.0₃0736"
$ws.Range("E18").Value = "  +1.02%  "
$ws.Range("D19").Value = "'223.26"
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +8.73%  "
$ws.Range("E23").Value = "  +4.00%  "
$ws.Range("D24").Value = "'9.29"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("D25").Value = "'146.98"
$ws.Range("E25").Value = "  +0.96%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +4.90%  "
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").Value = "'0.0516"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "1.263.54"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("D39").Value = "'0.841"
$ws.Range("E39").Value = "  +3.86%  "
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("E41").Value = "  +2.09%  "
$ws.Range("D42").Value = "'5.39"
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "1.801.27"
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("E44").Value = "  -4.09%  "
$ws.Range("D45").Value = "'62.00"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").Value = "'92.16"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("D47").Value = "'1.62"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  +2.08%  "
$ws.Range("D50").Value = "'7.67"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("E51").Value = "  +0.25%  "
